$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B10").Value = "Dump20160217-1"
$ws.Range("C10").Value = "Added new table (PRODUCT_CATEGORY)"

$ws.Range("C11").Select()
